$d = $word.ActiveDocument

$d.Content.Find.Execute("699×7=4893", $true, $false, $false, $false, $false, $true, 1, $false, "703×2=1406", 2)
$d.Content.Find.Execute("164×4=656", $true, $false, $false, $false, $false, $true, 1, $false, "908×4=3632", 2)
$d.Content.Find.Execute("481×2=962", $true, $false, $false, $false, $false, $true, 1, $false, "873×7=6111", 2)
$d.Content.Find.Execute("145×9=1305", $true, $false, $false, $false, $false, $true, 1, $false, "888×4=3552", 2)
$d.Content.Find.Execute("805×4=3220", $true, $false, $false, $false, $false, $true, 1, $false, "338×3=1014", 2)
$d.Content.Find.Execute("113×8=904", $true, $false, $false, $false, $false, $true, 1, $false, "430×6=2580", 2)
$d.Content.Find.Execute("766×5=3830", $true, $false, $false, $false, $false, $true, 1, $false, "120×2=240", 2)
$d.Content.Find.Execute("873×5=4365", $true, $false, $false, $false, $false, $true, 1, $false, "795×5=3975", 2)
$d.Content.Find.Execute("566×5=2830", $true, $false, $false, $false, $false, $true, 1, $false, "539×3=1617", 2)
$d.Content.Find.Execute("984×2=1968", $true, $false, $false, $false, $false, $true, 1, $false, "101×8=808", 2)
$d.Content.Find.Execute("197×7=1379", $true, $false, $false, $false, $false, $true, 1, $false, "923×3=2769", 2)
$d.Content.Find.Execute("402×4=1608", $true, $false, $false, $false, $false, $true, 1, $false, "609×3=1827", 2)
$d.Content.Find.Execute("701×9=6309", $true, $false, $false, $false, $false, $true, 1, $false, "265×6=1590", 2)
$d.Content.Find.Execute("436×5=2180", $true, $false, $false, $false, $false, $true, 1, $false, "938×3=2814", 2)
$d.Content.Find.Execute("198×4=792", $true, $false, $false, $false, $false, $true, 1, $false, "560×3=1680", 2)
$d.Content.Find.Execute("604×4=2416", $true, $false, $false, $false, $false, $true, 1, $false, "129×5=645", 2)
$d.Content.Find.Execute("997×5=4985", $true, $false, $false, $false, $false, $true, 1, $false, "509×4=2036", 2)
$d.Content.Find.Execute("570×3=1710", $true, $false, $false, $false, $false, $true, 1, $false, "549×4=2196", 2)
$d.Content.Find.Execute("699×9=6291", $true, $false, $false, $false, $false, $true, 1, $false, "481×8=3848", 2)
$d.Content.Find.Execute("540×8=4320", $true, $false, $false, $false, $false, $true, 1, $false, "356×6=2136", 2)
$d.Content.Find.Execute("410×8=3280", $true, $false, $false, $false, $false, $true, 1, $false, "290×3=870", 2)
$d.Content.Find.Execute("269×3=807", $true, $false, $false, $false, $false, $true, 1, $false, "456×3=1368", 2)
$d.Content.Find.Execute("933×8=7464", $true, $false, $false, $false, $false, $true, 1, $false, "394×7=2758", 2)
$d.Content.Find.Execute("517×3=1551", $true, $false, $false, $false, $false, $true, 1, $false, "196×2=392", 2)
$d.Content.Find.Execute("389×4=1556", $true, $false, $false, $false, $false, $true, 1, $false, "721×3=2163", 2)
